$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update calibration values (Legs Update)
$ws.Range("D5").Value = 1580
$ws.Range("G5").Value = 2050
$ws.Range("D8").Value = 1250

# Update selection on the sheet (Sesi 2 Update)
$ws.Range("F14").Select()
